$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 2.794415499126957
$ws.Range("B4").Value = 1.929662451756009
$ws.Range("C4").Value = 1.79986265093306
$ws.Range("D4").Value = -0.8854380943849716
$ws.Range("C5").Value = -0.3103476474035083
$ws.Range("D5").Value = -2.218321983005667
$ws.Range("E5").Value = 1.532898100704427
$ws.Range("F5").Value = -0.03788152406275502
$ws.Range("E6").Value = 1.0892423430376
$ws.Range("F6").Value = 0.3618709043640589
$ws.Range("G6").Value = -1.675983833549222
$ws.Range("H6").Value = 2.15807511757542
$ws.Range("G7").Value = 1.277667191469334
$ws.Range("H7").Value = 2.627768965131905
$ws.Range("I7").Value = 2.659314723144024
$ws.Range("J7").Value = 2.814292328656265
$ws.Range("I8").Value = 1.423575078814565
$ws.Range("J8").Value = 1.277551254953391
$ws.Range("K8").Value = 2.906157307553836
$ws.Range("L8").Value = 2.942581135514977
$ws.Range("K9").Value = 4.038794034641202
$ws.Range("L9").Value = 4.110668188518263
$ws.Range("M9").Value = 3.586987532670949
$ws.Range("N9").Value = 3.878230798954285
$ws.Range("M10").Value = 5.578558913710663
$ws.Range("N10").Value = 5.555223160690259
$ws.Range("O10").Value = 3.42596297413984
$ws.Range("P10").Value = 3.075158037444581
$ws.Range("O11").Value = 2.76275821580223
$ws.Range("P11").Value = 2.367041597905817
$ws.Range("Q11").Value = 1.999626938280241
$ws.Range("Q12").Value = 1.512800301290995
$ws.Range("R12").Value = 2.590730081186199
$ws.Range("S12").Value = 2.793289702145763
$ws.Range("R13").Value = 2.494417544901628
$ws.Range("S13").Value = 2.207308935472674
$ws.Range("T13").Value = 2.821516951149361
$ws.Range("U13").Value = 3.016819787229474
$ws.Range("T14").Value = 3.121124374623663
$ws.Range("U14").Value = 3.300883241600383
$ws.Range("V14").Value = 2.535635243126988
$ws.Range("W14").Value = 2.637488927515808
$ws.Range("X14").Value = 2.841516658941856
$ws.Range("V15").Value = 2.905506582474837
$ws.Range("W15").Value = 3.099927982210238
$ws.Range("X15").Value = 3.177336867742331
$ws.Range("Y15").Value = 2.46048248889319
$ws.Range("Z15").Value = 2.220932789361152
$ws.Range("AA15").Value = 2.167670286234991
$ws.Range("AB15").Value = 2.207173254521999
$ws.Range("Y16").Value = 2.574315362377289
$ws.Range("Z16").Value = 2.562407432124303
$ws.Range("AA16").Value = 2.556000311085604
$ws.Range("AB16").Value = 2.525788519949024
$ws.Range("AC16").Value = 1.496881353009161
$ws.Range("AD16").Value = 1.694296813984009
$ws.Range("AE16").Value = 1.629686186121027
$ws.Range("AF16").Value = 1.619782579158202
$ws.Range("AC17").Value = 1.83915572102098
$ws.Range("AD17").Value = 1.811820461872138
$ws.Range("AE17").Value = 1.633300070291677
$ws.Range("AF17").Value = 0.3071129274195616
$ws.Range("AG17").Value = 2.257237375640031
$ws.Range("AH17").Value = 2.059057776028594
$ws.Range("AI17").Value = 1.937440955395164
$ws.Range("AJ17").Value = 1.388548717051186
$ws.Range("AG18").Value = 2.213620378726788
$ws.Range("AH18").Value = 2.610416778758373
$ws.Range("AI18").Value = 2.179993174715689
$ws.Range("AJ18").Value = 0.890977499942136
$ws.Range("AK18").Value = 2.058050235820175
$ws.Range("AL18").Value = 2.29686889447267
$ws.Range("AM18").Value = 3.014134262744617
$ws.Range("AN18").Value = 1.437806261771213
$ws.Range("AK19").Value = 2.849992723907335
$ws.Range("AL19").Value = 3.479464952554112
$ws.Range("AM19").Value = 4.125217580302332
$ws.Range("AN19").Value = 2.148304186541194
$ws.Range("AO19").Value = 2.599913004672616
$ws.Range("AP19").Value = 2.337862417976333
$ws.Range("AQ19").Value = 2.577692526489739
$ws.Range("AR19").Value = -1.390622874876313
$ws.Range("AO20").Value = 1.343460690969822
$ws.Range("AP20").Value = 0.8329290289207147
$ws.Range("AQ20").Value = 1.539533176834884
$ws.Range("AR20").Value = -2.093034802586002
$ws.Range("AS20").Value = 0.5029237023806754
$ws.Range("AT20").Value = 0.4565833992175916
$ws.Range("AU20").Value = 1.187924830910969
$ws.Range("AV20").Value = 0.8461784325530575
$ws.Range("AS21").Value = -0.7352716516441982
$ws.Range("AT21").Value = -0.7869205535448565
$ws.Range("AU21").Value = -0.3489999547360179
$ws.Range("AV21").Value = -0.9869022883377543
$ws.Range("AW21").Value = 1.005480064500386
$ws.Range("AX21").Value = 0.5434772144153888
$ws.Range("AY21").Value = 0.5503752294844233
$ws.Range("AZ21").Value = 1.903706680019468
$ws.Range("AW22").Value = 0.5134406156019233
$ws.Range("AX22").Value = 0.03073175363270675
$ws.Range("AY22").Value = -0.04341979710753563
$ws.Range("AZ22").Value = 0.6808804886353492
$ws.Range("BA22").Value = 2.403526819519342
$ws.Range("BA23").Value = 1.910815645575914
